$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 383.75
$ws.Range("I28").Value = 121.84615
$ws.Range("K28").Value = 121.84615
$ws.Range("M28").Value = 363.15385
$ws.Range("H101").Value = 4386.1763
$ws.Range("I101").Value = 1324.909
$ws.Range("K101").Value = 3974.727
$ws.Range("M101").Value = -2352.727
$ws.Range("H111").Value = 680.38464
$ws.Range("I111").Value = 807.4
$ws.Range("J111").Value = 601
$ws.Range("K111").Value = 2422.2
$ws.Range("L111").Value = 1803
$ws.Range("M111").Value = 644.8000000000002
$ws.Range("N111").Value = -7937
$ws.Range("H129").Value = 3206534.8
$ws.Range("I129").Value = 683.625
$ws.Range("J129").Value = 8335897
$ws.Range("K129").Value = 2050.875
$ws.Range("L129").Value = 25007691
$ws.Range("M129").Value = 2949.125
$ws.Range("N129").Value = -25017691
$ws.Range("H131").Value = 12244
$ws.Range("I131").Value = 9085.625
$ws.Range("J131").Value = 20666.334
$ws.Range("K131").Value = 27256.875
$ws.Range("L131").Value = 61999.00199999999
$ws.Range("M131").Value = -22216.875
$ws.Range("N131").Value = -72079.00199999999
$ws.Range("H132").Value = 1974.3182
$ws.Range("I132").Value = 1974.3182
$ws.Range("K132").Value = 5922.9546
$ws.Range("M132").Value = -3392.9546
$ws.Range("H138").Value = 3994.9473
$ws.Range("I138").Value = 6210.8887
$ws.Range("K138").Value = 18632.6661
$ws.Range("M138").Value = -13492.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1035.8667
$ws.Range("J2").Value = 2145.2
$ws.Range("L2").Value = 2145.2
$ws.Range("N2").Value = -2371.2
$ws.Range("H14").Value = 500
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H45").Value = 6183.3335
$ws.Range("I45").Value = 4885.7144
$ws.Range("K45").Value = 4885.7144
$ws.Range("M45").Value = -4508.7144
$ws.Range("H61").Value = 5458.091
$ws.Range("I61").Value = 5458.091
$ws.Range("K61").Value = 5458.091
$ws.Range("M61").Value = -5246.091
$ws.Range("H116").Value = 1035.8667
$ws.Range("J116").Value = 2145.2
$ws.Range("L116").Value = 2145.2
$ws.Range("N116").Value = -6733.2
$ws.Range("H136").Value = 5458.091
$ws.Range("I136").Value = 5458.091
$ws.Range("K136").Value = 16374.273
$ws.Range("M136").Value = -13824.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1035.8667
$ws.Range("J3").Value = 2145.2
$ws.Range("L3").Value = 2145.2
$ws.Range("N3").Value = -2373.2
$ws.Range("H94").Value = 823.2727
$ws.Range("I94").Value = 906.2222
$ws.Range("K94").Value = 906.2222
$ws.Range("M94").Value = -455.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 189
$ws.Range("I19").Value = 243.83333
$ws.Range("J19").Value = 24.5
$ws.Range("K19").Value = 243.83333
$ws.Range("L19").Value = 24.5
$ws.Range("M19").Value = -73.83332999999999
$ws.Range("N19").Value = -364.5
$ws.Range("H24").Value = 189
$ws.Range("I24").Value = 243.83333
$ws.Range("J24").Value = 24.5
$ws.Range("K24").Value = 243.83333
$ws.Range("L24").Value = 24.5
$ws.Range("M24").Value = -73.83332999999999
$ws.Range("N24").Value = -364.5
$ws.Range("H58").Value = 7248.6787
$ws.Range("I58").Value = 2586.1177
$ws.Range("K58").Value = 2586.1177
$ws.Range("M58").Value = -2383.1177
$ws.Range("H97").Value = 25268.5
$ws.Range("J97").Value = 26025.428
$ws.Range("L97").Value = 26025.428
$ws.Range("N97").Value = -28007.428
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H136").Value = 7248.6787
$ws.Range("I136").Value = 2586.1177
$ws.Range("K136").Value = 7758.353099999999
$ws.Range("M136").Value = -5208.353099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 571788.7
$ws.Range("I7").Value = 333753.66
$ws.Range("J7").Value = 1999999
$ws.Range("K7").Value = 1001260.98
$ws.Range("L7").Value = 5999997
$ws.Range("M7").Value = -1001148.98
$ws.Range("N7").Value = -6000221
$ws.Range("H34").Value = 569.2857
$ws.Range("J34").Value = 497
$ws.Range("L34").Value = 1491
$ws.Range("N34").Value = -1659
$ws.Range("H107").Value = 417200.4
$ws.Range("I107").Value = 390.69232
$ws.Range("J107").Value = 909793.75
$ws.Range("K107").Value = 1172.07696
$ws.Range("L107").Value = 2729381.25
$ws.Range("M107").Value = 747.9230400000001
$ws.Range("N107").Value = -2733221.25
$ws.Range("H132").Value = 83334550
$ws.Range("I132").Value = 250000100
$ws.Range("K132").Value = 2250000900
$ws.Range("M132").Value = -2249998370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 500
$ws.Range("J3").Value = 500
$ws.Range("L3").Value = 500
$ws.Range("N3").Value = -732
$ws.Range("H5").Value = 50001
$ws.Range("J5").Value = 50001
$ws.Range("L5").Value = 50001
$ws.Range("N5").Value = -50225
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H102").Value = 3258.3572
$ws.Range("I102").Value = 2293.5
$ws.Range("K102").Value = 2293.5
$ws.Range("M102").Value = -671.5
$ws.Range("H113").Value = 670241.7
$ws.Range("I113").Value = 1334483.4
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 1334483.4
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -1332313.4
$ws.Range("N113").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2650
$ws.Range("J10").Value = 2650
$ws.Range("L10").Value = 2650
$ws.Range("N10").Value = -2930
$ws.Range("H40").Value = 2765.5833
$ws.Range("I40").Value = 2698.0715
$ws.Range("K40").Value = 2698.0715
$ws.Range("M40").Value = -2562.0715
$ws.Range("H61").Value = 63263.375
$ws.Range("I61").Value = 77769.30499999999
$ws.Range("K61").Value = 77769.30499999999
$ws.Range("M61").Value = -77567.30499999999
$ws.Range("H68").Value = 6297.25
$ws.Range("I68").Value = 2595
$ws.Range("K68").Value = 2595
$ws.Range("M68").Value = -1846
$ws.Range("H71").Value = 6297.25
$ws.Range("I71").Value = 2595
$ws.Range("K71").Value = 12975
$ws.Range("M71").Value = -9231
$ws.Range("H93").Value = 17349.691
$ws.Range("I93").Value = 2494.5557
$ws.Range("K93").Value = 2494.5557
$ws.Range("M93").Value = -1246.5557
$ws.Range("H113").Value = 63263.375
$ws.Range("I113").Value = 77769.30499999999
$ws.Range("K113").Value = 77769.30499999999
$ws.Range("M113").Value = -75599.30499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 13000
$ws.Range("I21").Value = 13000
$ws.Range("K21").Value = 13000
$ws.Range("M21").Value = -12765
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H35").Value = 13000
$ws.Range("I35").Value = 13000
$ws.Range("K35").Value = 13000
$ws.Range("M35").Value = -12710
$ws.Range("H51").Value = 22535
$ws.Range("I51").Value = 22535
$ws.Range("K51").Value = 22535
$ws.Range("M51").Value = -22025
$ws.Range("H99").Value = 40586.668
$ws.Range("J99").Value = 40586.668
$ws.Range("L99").Value = 40586.668
$ws.Range("N99").Value = -46576.668
$ws.Range("H122").Value = 3443.0613
$ws.Range("I122").Value = 2759.7896
$ws.Range("K122").Value = 8279.3688
$ws.Range("M122").Value = -5829.3688
$ws.Range("H132").Value = 4042.9092
$ws.Range("I132").Value = 3526.2104
$ws.Range("K132").Value = 10578.6312
$ws.Range("M132").Value = -8048.6312
$ws.Range("H136").Value = 3560.0312
$ws.Range("I136").Value = 1892.5834
$ws.Range("K136").Value = 5677.7502
$ws.Range("M136").Value = -3127.7502
